$d = $word.ActiveDocument
$r = $d.Content
$cur = $r.Font.Name
Write-Output "font=$cur"
